$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a cell's value as literal text, preventing Excel from
# auto-converting numeric-looking strings (e.g. "622.18") into actual
# numbers (which would lose/alter precision and the "t=inlineStr/s"
# string-cell representation). We temporarily force the Text number
# format, assign the value, then restore the cell's original style so
# no new style is left behind in the workbook.
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = $origStyle
}

Set-TextValue "D2" "92.411.59"
$ws.Range("E2").Value = "  -5.82%  "

Set-TextValue "D3" "3.325.57"
$ws.Range("E3").Value = "  -4.49%  "

$ws.Range("E4").Value = "  -0.17%  "

Set-TextValue "D5" "228.37"
$ws.Range("E5").Value = "  -9.24%  "

Set-TextValue "D6" "622.18"
$ws.Range("E6").Value = "  -6.33%  "

Set-TextValue "D7" "1.33"
$ws.Range("E7").Value = "  -9.81%  "

Set-TextValue "D8" "0.377"
$ws.Range("E8").Value = "  -11.57%  "

Set-TextValue "D9" "0.999"
$ws.Range("E9").Value = "  -0.06%  "

Set-TextValue "D10" "0.904"
$ws.Range("E10").Value = "  -14.48%  "

Set-TextValue "D11" "3.323.01"
$ws.Range("E11").Value = "  -4.48%  "

$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D12" "40.56"
$ws.Range("E12").Value = "  -11.07%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D13" "0.189"
$ws.Range("E13").Value = "  -9.94%  "

Set-TextValue "D14" "92.041.85"
$ws.Range("E14").Value = "  -6.07%  "

Set-TextValue "D15" "5.83"
$ws.Range("E15").Value = "  -6.09%  "

Set-TextValue "D16" "3.936.98"
$ws.Range("E16").Value = "  -4.56%  "

Set-TextValue "D17" "0.0000240"
$ws.Range("E17").Value = "  -7.67%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D18" "7.89"
$ws.Range("E18").Value = "  -12.05%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D19" "3.317.67"
$ws.Range("E19").Value = "  -4.83%  "

Set-TextValue "D20" "16.66"
$ws.Range("E20").Value = "  -10.76%  "

Set-TextValue "D21" "10.77"
$ws.Range("E21").Value = "  -8.97%  "

Set-TextValue "D22" "484.79"
$ws.Range("E22").Value = "  -7.30%  "

Set-TextValue "D23" "3.16"
$ws.Range("E23").Value = "  -7.23%  "

Set-TextValue "D24" "0.441"
$ws.Range("E24").Value = "  -14.42%  "

Set-TextValue "D25" "0.0000179"
$ws.Range("E25").Value = "  -11.27%  "

Set-TextValue "D26" "6.04"
$ws.Range("E26").Value = "  -10.84%  "

Set-TextValue "D27" "88.93"
$ws.Range("E27").Value = "  -8.86%  "

$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D28" "11.31"
$ws.Range("E28").Value = "  -10.99%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D29" "1.00"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D30" "11.04"
$ws.Range("E30").Value = "  -11.32%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D31" "2.58"
$ws.Range("E31").Value = "  -9.26%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D32" "0.131"
$ws.Range("E32").Value = "  -9.44%  "

$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D33" "0.991"
$ws.Range("E33").Value = "  -0.73%  "

$ws.Range("B34").Value = "Cronos"
$ws.Range("C34").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D34" "0.168"
$ws.Range("E34").Value = "  -11.61%  "

$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D35" "28.13"
$ws.Range("E35").Value = "  -9.29%  "

$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D36" "0.515"
$ws.Range("E36").Value = "  -13.75%  "

$ws.Range("B37").Value = "USDe"
$ws.Range("C37").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D37" "1.00"
$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D38" "519.52"
$ws.Range("E38").Value = "  -0.74%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D39" "7.28"
$ws.Range("E39").Value = "  -7.64%  "

$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D40" "1.35"
$ws.Range("E40").Value = "  -9.63%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D41" "0.144"
$ws.Range("E41").Value = "  -7.52%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D42" "0.859"
$ws.Range("E42").Value = "  -5.74%  "

Set-TextValue "D43" "23.97"
$ws.Range("E43").Value = "  -1.81%  "

$ws.Range("B44").Value = "MantraDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue "D44" "3.56"
$ws.Range("E44").Value = "  -2.19%  "

$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D45" "1.65"
$ws.Range("E45").Value = "  -5.18%  "

$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D46" "5.30"
$ws.Range("E46").Value = "  -7.34%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D47" "2.12"
$ws.Range("E47").Value = "  -4.47%  "

$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D48" "52.29"
$ws.Range("E48").Value = "  -5.09%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D49" "0.0387"
$ws.Range("E49").Value = "  -10.15%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D50" "7.75"
$ws.Range("E50").Value = "  -10.52%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D51" "3.01"
$ws.Range("E51").Value = "  -7.30%  "

Write-Output "Done applying cryptos update"
